# edit.ps1
# Applies the "20 oct 17pm snapshot" change to the Kmeans snapshot log table:
#   1. Removes the stray <w:rFonts w:hint="cs"/> paragraph-mark hint from the
#      three paragraphs of the last existing row (17.10 / 15:00 / the "COMPARE"
#      cell's second paragraph) that no longer need it.
#   2. Appends a new table row (20.10 / 17:00 / bug-fix note) at the end of the
#      log table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-ParagraphXml($paragraph, [string]$innerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($pkg)
}

# --- Step 1: drop the paragraph-mark rFonts hint from the last row's cells ---
$lastRowIndex = $t.Rows.Count

$cell1 = $t.Cell($lastRowIndex, 1)
Set-ParagraphXml $cell1.Range.Paragraphs.Item(1) '<w:p><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>17.10</w:t></w:r></w:p>'

$cell2 = $t.Cell($lastRowIndex, 2)
Set-ParagraphXml $cell2.Range.Paragraphs.Item(1) '<w:p><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>15:00</w:t></w:r></w:p>'

$cell3 = $t.Cell($lastRowIndex, 3)
Set-ParagraphXml $cell3.Range.Paragraphs.Item(2) '<w:p><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>יש תקלה עם שמירה של 0 בין טרנסקציות של חישוב שלם שצריך לתקן ב</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/></w:rPr><w:t>S</w:t></w:r><w:r><w:t xml:space="preserve">COREBOARD </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> או במוניטור שמעביר לו על החוט.</w:t></w:r></w:p>'

# --- Step 2: append the new "20.10 / 17:00" row ---
$newRow = $t.Rows.Add()
$newRow.Height = 554
$newRowIndex = $newRow.Index

$newCell1 = $t.Cell($newRowIndex, 1)
Set-ParagraphXml $newCell1.Range.Paragraphs.Item(1) '<w:p><w:pPr><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>20.10</w:t></w:r></w:p>'

$newCell2 = $t.Cell($newRowIndex, 2)
Set-ParagraphXml $newCell2.Range.Paragraphs.Item(1) '<w:p><w:pPr><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>17:00</w:t></w:r></w:p>'

$newCell3 = $t.Cell($newRowIndex, 3)
Set-ParagraphXml $newCell3.Range.Paragraphs.Item(1) '<w:p><w:pPr><w:rPr><w:rFonts w:hint="cs"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>תיקון באגים של שגיאה מצטברת, כרגע יש שגיאה אקראית כתלות בהגרלה, סנטרואיד 8 תמיד שגיאה אפסית, כל השאר משתנה, המקסימום זה 5000 אלפיות שזה לא מעט..</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> עבור סנטרואיד כלשהו.</w:t></w:r></w:p>'

Write-Host "done"
